$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text fixes ---
$ws.Range("A2").Value = "make span: 127"
$ws.Range("A3").Value = "Machines"

# --- Numeric data row 4 (per-machine production-time totals) ---
$ws.Range("A4").Value = 33
$ws.Range("B4").Value = 48
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 6

# --- Row 5 stays the same (machine indices) ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4

# --- Order labels rows 6-8 (re-sorted order list after fixing the
#     "added production time twice" bug) ---
$ws.Range("A6").Value = "Order 4 - 20"
$ws.Range("B6").Value = "Order 5 - 40"
$ws.Range("C6").Value = "Order 7 - 57"
$ws.Range("D6").Value = "Order 6 - 26"

$ws.Range("A7").Value = "Order 9 - 68"
$ws.Range("B7").Value = "Order 2 - 79"
$ws.Range("C7").Value = "Order 10 - 127"
$ws.Range("D7").Value = "Order 1 - 87"

$ws.Range("A8").Value = "Order 3 - 94"
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = "Order 8 - 121"

# --- Column widths: only column A keeps a custom width now (~16.14 chars);
#     B and C revert to the workbook's standard/default width.
#     (15.25 is the input that the host's pixel-quantized ColumnWidth
#     rounds to the value closest to the target 16.140625.) ---
$ws.Columns.Item(1).ColumnWidth = 15.25
$ws.Columns.Item(2).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(3).ColumnWidth = $ws.StandardWidth
